$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Haste")

# Update row 3 values (Braunschweig (RE) row) - text stays, numbers change
$ws.Range("B3").Value() = 18
$ws.Range("C3").Value() = 34
$ws.Range("D3").Value() = 28
$ws.Range("E3").Value() = 29
$ws.Range("F3").Value() = 2

# Delete old rows 7 and 8 (Nienburg (S) and Hannover Hbf (RE))
$ws.Rows.Item(7).Delete() | Out-Null
$ws.Rows.Item(7).Delete() | Out-Null

# Update row 6: replace "Hannover Hbf (S2)" with "Weetzen (S)" and new numbers
$ws.Range("A6").Value() = "Weetzen (S)"
$ws.Range("B6").Value() = 31
$ws.Range("C6").Value() = 34
$ws.Range("D6").Value() = 26
$ws.Range("E6").Value() = 27
$ws.Range("F6").Value() = 3

$ws.Rows.Item(6).AutoFit() | Out-Null

$ws.Range("F7").Select() | Out-Null

Write-Host "done"
